$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("C2").Value = 0.5785953177257525
$ws.Range("D2").Value = 0.9719101123595506
$ws.Range("E2").Value = 0.7253668763102725
$ws.Range("F2").Value = 0.8555885262116716
$ws.Range("G2").Value = 0.9471467677405769
$ws.Range("H2").Value = 0.7777356955491029
$ws.Range("I2").Value = 519
$ws.Range("J2").Value = 378
$ws.Range("K2").Value = 156
$ws.Range("L2").Value = 15

# ---- Classification Report sheet ----
$ws = $wb.Worksheets.Item("Classification Report")
$ws.Range("B2").Value = 0.9122807017543859
$ws.Range("C2").Value = 0.2921348314606741
$ws.Range("D2").Value = 0.4425531914893617

$ws.Range("B3").Value = 0.5785953177257525
$ws.Range("C3").Value = 0.9719101123595506
$ws.Range("D3").Value = 0.7253668763102725

$ws.Range("B5").Value = 0.7454380097400692
$ws.Range("C5").Value = 0.6320224719101124
$ws.Range("D5").Value = 0.5839600338998171

$ws.Range("B6").Value = 0.7454380097400691
$ws.Range("D6").Value = 0.5839600338998171

# ---- Confusion Matrix sheet ----
$ws = $wb.Worksheets.Item("Confusion Matrix")
$ws.Range("B2").Value = 156
$ws.Range("C2").Value = 378
$ws.Range("B3").Value = 15
$ws.Range("C3").Value = 519
